$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.1960132890365449
$ws.Range("C2").Value = 0.5481727574750831
$ws.Range("J2").Value = 0.026578073089701
$ws.Range("P2").Value = 0.1162790697674419
$ws.Range("S2").Value = 0.1129568106312292
$ws.Range("C3").Value = 0.01775147928994083
$ws.Range("J3").Value = 0.01775147928994083
$ws.Range("P3").Value = 0.6686390532544378
$ws.Range("S3").Value = 0.2958579881656805
$ws.Range("P4").Value = 0.7692307692307693
$ws.Range("S4").Value = 0.2307692307692308
$ws.Range("B6").Value = 0.08900523560209424
$ws.Range("D6").Value = 0.01570680628272251
$ws.Range("F6").Value = 0.04712041884816754
$ws.Range("J6").Value = 0.225130890052356
$ws.Range("O6").Value = 0.005235602094240838
$ws.Range("Q6").Value = 0.1675392670157068
$ws.Range("R6").Value = 0.05759162303664921
$ws.Range("S6").Value = 0.3926701570680629
$ws.Range("B7").Value = 0.1472868217054264
$ws.Range("D7").Value = 0.007751937984496124
$ws.Range("E7").Value = 0.007751937984496124
$ws.Range("F7").Value = 0.007751937984496124
$ws.Range("J7").Value = 0.1472868217054264
$ws.Range("O7").Value = 0.007751937984496124
$ws.Range("Q7").Value = 0.248062015503876
$ws.Range("R7").Value = 0.06201550387596899
$ws.Range("S7").Value = 0.3643410852713178
$ws.Range("B8").Value = 0.09978308026030369
$ws.Range("D8").Value = 0.01518438177874186
$ws.Range("F8").Value = 0.05206073752711497
$ws.Range("J8").Value = 0.1106290672451193
$ws.Range("O8").Value = 0.01301518438177874
$ws.Range("Q8").Value = 0.2212581344902386
$ws.Range("R8").Value = 0.0737527114967462
$ws.Range("S8").Value = 0.4143167028199566
$ws.Range("B9").Value = 0.09883720930232558
$ws.Range("D9").Value = 0.005813953488372093
$ws.Range("F9").Value = 0.06395348837209303
$ws.Range("J9").Value = 0.1104651162790698
$ws.Range("O9").Value = 0.005813953488372093
$ws.Range("Q9").Value = 0.2151162790697674
$ws.Range("R9").Value = 0.05232558139534884
$ws.Range("S9").Value = 0.4476744186046512
$ws.Range("B10").Value = 0.0982532751091703
$ws.Range("D10").Value = 0.01965065502183406
$ws.Range("E10").Value = 0.000727802037845706
$ws.Range("F10").Value = 0.06186317321688501
$ws.Range("J10").Value = 0.1324599708879185
$ws.Range("O10").Value = 0.008005822416302766
$ws.Range("Q10").Value = 0.2634643377001455
$ws.Range("R10").Value = 0.08296943231441048
$ws.Range("S10").Value = 0.3326055312954876
$ws.Range("F11").Value = 0.004784688995215311
$ws.Range("G11").Value = 0.1483253588516746
$ws.Range("J11").Value = 0.08133971291866028
$ws.Range("K11").Value = 0.2200956937799043
$ws.Range("L11").Value = 0.5167464114832536
$ws.Range("S11").Value = 0.02870813397129187
$ws.Range("G12").Value = 0.6752136752136753
$ws.Range("J12").Value = 0.2564102564102564
$ws.Range("K12").Value = 0.008547008547008548
$ws.Range("L12").Value = 0.05982905982905983
$ws.Range("G13").Value = 0.5833333333333334
$ws.Range("S13").Value = 0.08333333333333333
$ws.Range("G14").Value = 0.6666666666666666
$ws.Range("J14").Value = 0.3333333333333333
$ws.Range("H15").Value = 0.1208791208791209
$ws.Range("I15").Value = 0.08791208791208792
$ws.Range("J15").Value = 0.5054945054945055
$ws.Range("K15").Value = 0.06593406593406594
$ws.Range("M15").Value = 0.02197802197802198
$ws.Range("O15").Value = 0.05494505494505494
$ws.Range("S15").Value = 0.1428571428571428
$ws.Range("F16").Value = 0.02923976608187134
$ws.Range("H16").Value = 0.1637426900584795
$ws.Range("I16").Value = 0.08187134502923976
$ws.Range("J16").Value = 0.4502923976608187
$ws.Range("K16").Value = 0.08187134502923976
$ws.Range("M16").Value = 0.005847953216374269
$ws.Range("O16").Value = 0.02923976608187134
$ws.Range("S16").Value = 0.1578947368421053
$ws.Range("F17").Value = 0.01071428571428571
$ws.Range("H17").Value = 0.2017857142857143
$ws.Range("I17").Value = 0.07321428571428572
$ws.Range("J17").Value = 0.4517857142857143
$ws.Range("K17").Value = 0.05892857142857143
$ws.Range("M17").Value = 0.01607142857142857
$ws.Range("N17").Value = 0.001785714285714286
$ws.Range("O17").Value = 0.04642857142857143
$ws.Range("S17").Value = 0.1392857142857143
$ws.Range("F18").Value = 0.01149425287356322
$ws.Range("H18").Value = 0.1666666666666667
$ws.Range("I18").Value = 0.05747126436781609
$ws.Range("J18").Value = 0.4770114942528735
$ws.Range("K18").Value = 0.07471264367816093
$ws.Range("M18").Value = 0.005747126436781609
$ws.Range("O18").Value = 0.08045977011494253
$ws.Range("S18").Value = 0.1264367816091954
$ws.Range("F19").Value = 0.01535974130962005
$ws.Range("H19").Value = 0.2158447857720291
$ws.Range("I19").Value = 0.07356507679870655
$ws.Range("J19").Value = 0.4066289409862571
$ws.Range("K19").Value = 0.06952303961196443
$ws.Range("M19").Value = 0.01535974130962005
$ws.Range("O19").Value = 0.07194826192400972
$ws.Range("S19").Value = 0.1285367825383993
